$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 19; $r -le 50; $r++) {
    $ws.Cells.Item($r, 9).Value = 4.435333333333333
}
